$d = $word.ActiveDocument

# --- Remove the trailing space run after "Working on UI/UX" ---
$lastPara = $d.Paragraphs.Last
$r = $lastPara.Range
$trimEnd = $r.End - 1          # position right before the paragraph mark
$trimStart = $trimEnd - 1      # position right before the trailing space
$spaceRange = $d.Range($trimStart, $trimEnd)
if ($spaceRange.Text -eq " ") {
    $spaceRange.Text = ""
}

# --- Insert a fresh empty paragraph to use as the anchor for the new content.
#     It inherits the "Working on UI/UX" paragraph's ListParagraph/numId=7
#     formatting, which InsertXML below will push down onto the final
#     (still-empty) list item, exactly matching the target structure. ---
$last = $d.Paragraphs.Last
$r = $last.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$anchor = $d.Paragraphs.Last
$rng = $anchor.Range
$rng.Collapse(1)

$xml = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:r><w:t>January 10</w:t></w:r><w:r><w:rPr><w:vertAlign w:val="superscript"/></w:rPr><w:t>th</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr></w:pPr><w:r><w:t>Finish UI/UX updates</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Update </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>github</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> with new work</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData></pkg:part>
</pkg:package>
"@

$rng.InsertXML($xml)
